$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (ciudad, tecnico, modo, costo_aprox_fase1_uf, dias_aprox_fase1)
$data = @(
    @("Viña del Mar", "Luis", "terrestre", 43.6863, 9),
    @("Talca", "Luis", "terrestre", 16.0235, 3),
    @("Antofagasta", "Orlando", "terrestre", 19.82393333333333, 4),
    @("Calama", "Orlando", "terrestre", 7.286666666666667, 2),
    @("Iquique", "Orlando", "terrestre", 10.33776666666667, 2),
    @("Temuco", "Jimmy", "terrestre", 25.81096666666667, 5),
    @("Concepcion", "Jimmy", "terrestre", 28.2453, 6),
    @("Chillan", "Jimmy", "terrestre", 9.950000000000001, 3),
    @("Puerto Montt", "Jimmy", "terrestre", 15.81144, 3),
    @("Osorno", "Jimmy", "terrestre", 16.268, 3),
    @("Los Angeles", "Jimmy", "terrestre", 14.5658, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
